$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Activate()

$ws.Range("D3").Value = "new changes"

$ws.Range("D4").Select()
